# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have been generated: status text, handback
# file/datetime columns are populated, links are added, and a few
# columns are widened so the new content is readable.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$mdFileName       = "41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f.md"
$mdUrl            = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a40e2c853ad572fcfc300acd0f6da4de29d4b9c/e2e/41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f.md"

# Widths below are expressed so that, after Excel's internal column-width
# quantization (stored width = ColumnWidth + 0.8333333333333334, rounded
# to the nearest 1/6), they land on the intended rendered widths.
$wideWidth   = 29.144371396019366   # -> stored width ~29.98 (was ~17.22)
$fortyWidth  = 39.166666666666664   # -> stored width 40 (was ~18.65 / ~21.71)

# ---------------------------------------------------------------------
# Overview sheet: status columns for zh-cn (E2) and de-de (F2), and make
# those columns wider to fit the longer status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E1").EntireColumn.ColumnWidth = $wideWidth
$overview.Range("F1").EntireColumn.ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C1").EntireColumn.ColumnWidth = $wideWidth
$zhcn.Range("I1").EntireColumn.ColumnWidth = $fortyWidth
$zhcn.Range("J1").EntireColumn.ColumnWidth = $fortyWidth

# Latest Target File (I2) becomes a hyperlink to the handed-back markdown
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFileName)
# Latest Handback File (J2) - same xliff file name used for the handoff
$zhcn.Range("J2").Value = "41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f.c4004000eca71d9d6295042713ef2b07d8932e65.zh-cn.xlf"
# Latest Handback DateTime (K2)
$zhcn.Range("K2").Value = "2016-09-02 11:09:31"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C1").EntireColumn.ColumnWidth = $wideWidth
$dede.Range("I1").EntireColumn.ColumnWidth = $fortyWidth
$dede.Range("J1").EntireColumn.ColumnWidth = $fortyWidth

# Latest Target File (I2) becomes a hyperlink to the handed-back markdown
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFileName)
# Latest Handback File (J2) - same xliff file name used for the handoff
$dede.Range("J2").Value = "41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f.c4004000eca71d9d6295042713ef2b07d8932e65.de-de.xlf"
# Latest Handback DateTime (K2)
$dede.Range("K2").Value = "2016-09-02 11:09:39"
